$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44533
$ws.Range("M2").Value = 100

# Row 3
$ws.Range("D3").Value = 44357
$ws.Range("N3").Value = 14000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 14500
$ws.Range("S3").Value = 725

# Row 4
$ws.Range("D4").Value = 44761
$ws.Range("M4").Value = 100
$ws.Range("O4").Value = 21000
$ws.Range("P4").Value = 20500
$ws.Range("S4").Value = 1025

# Row 5
$ws.Range("D5").Value = 44890
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 23000
$ws.Range("P5").Value = 22250
$ws.Range("S5").Value = 1112

# Row 6
$ws.Range("D6").Value = 44893
$ws.Range("M6").Value = 80
$ws.Range("N6").Value = 21000
$ws.Range("O6").Value = 22000
$ws.Range("P6").Value = 21625
$ws.Range("S6").Value = 1081

# Row 8
$ws.Range("D8").Value = 44320
$ws.Range("M8").Value = 80

# Row 9
$ws.Range("D9").Value = 44792
$ws.Range("M9").Value = 100

# Row 10
$ws.Range("D10").Value = 44798
$ws.Range("M10").Value = 80
$ws.Range("N10").Value = 21000
$ws.Range("O10").Value = 22000
$ws.Range("P10").Value = 21500
$ws.Range("S10").Value = 1075
